$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 90 - this shifts existing rows 90..183 down to 91..184
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly data record
$ws.Cells.Item(90, 1).Value = 8
$ws.Cells.Item(90, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(90, 3).Value = "Coquimbo"
$ws.Cells.Item(90, 4).Value = 44601
$ws.Cells.Item(90, 5).Value = 4
$ws.Cells.Item(90, 6).Value = 100112031
$ws.Cells.Item(90, 7).Value = "Poroto verde"
$ws.Cells.Item(90, 8).Value = "Magnum"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 360
$ws.Cells.Item(90, 11).Value = 37000
$ws.Cells.Item(90, 12).Value = 38000
$ws.Cells.Item(90, 13).Value = 37500
$ws.Cells.Item(90, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(90, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(90, 16).Value = 1500
$ws.Cells.Item(90, 17).Value = 25
$ws.Cells.Item(90, 18).Value = "Hortaliza"
